# Real Madrid fixture update: remove the Celta de Vigo match (row 1),
# shift remaining rows up, refresh the "tickets remaining" counts, and
# append the newly scheduled Atlético de Madrid match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (shared string), never letting
# Excel's automatic number/date detection turn it into a numeric cell.
# We do this by building the text via a formula (="literal"), copying
# that single cell, and pasting-special (values only) into the target -
# this yields a plain text cell without adding/touching any cell styles.
$scratch = $ws.Range("ZZ1")
function Set-TextValue {
    param($range, [string]$text)
    $escaped = $text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $range.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
    $scratch.ClearContents()
}

# 1) Delete the first match (Real Madrid CF - Celta de Vigo) entirely;
#    this shifts every following row up by one.
$ws.Rows.Item(1).Delete()

# 2) Refresh the ticket counts in column C for the shifted rows.
Set-TextValue $ws.Range("C1") "53"
Set-TextValue $ws.Range("C2") "51"
Set-TextValue $ws.Range("C3") "56"
Set-TextValue $ws.Range("C4") "36"
Set-TextValue $ws.Range("C5") "57"
Set-TextValue $ws.Range("C6") "36"
Set-TextValue $ws.Range("C7") "59"
Set-TextValue $ws.Range("C8") "37"
Set-TextValue $ws.Range("C9") "36"

# 3) Append the new fixture row for the Atlético de Madrid match.
$ws.Range("A10").Value = "Real Madrid CF - Atlético de Madrid"
Set-TextValue $ws.Range("B10") "12/12/2021"
Set-TextValue $ws.Range("C10") "54"

$scratch.ClearContents()
